# Regenerate the "K" (strikeouts) column (column G) in the save_data sheet
# for giolito_lucas.xlsx. The values previously reflected a different stat
# (Strike#) and are being recalculated/rewritten to hold strikeout counts (K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 38 (row r corresponds to start index r-2)
$kValues = @{
    2  = 4
    3  = 5
    4  = 6
    5  = 12
    6  = 9
    7  = 3
    8  = 6
    9  = 9
    10 = 5
    11 = 7
    12 = 2
    13 = 5
    14 = 9
    15 = 5
    16 = 5
    17 = 4
    18 = 9
    19 = 10
    20 = 5
    21 = 8
    22 = 10
    23 = 5
    24 = 4
    25 = 3
    26 = 5
    27 = 9
    28 = 7
    29 = 6
    30 = 5
    31 = 7
    32 = 7
    33 = 3
    34 = 6
    35 = 6
    36 = 6
    37 = 6
    38 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
